$wb = $excel.ActiveWorkbook

# --- Update conversion text on sheet "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.36 = 12926.21 pesos`n✅ 12926.21 pesos = 3.34 = 962.45 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update tasas values on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 298
$wsTasas.Range("O10").Value = 3852.01
$wsTasas.Range("N12").Value = 3868
$wsTasas.Range("O12").Value = 288
